$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Golf")

# Update header labels in row 160 to drop the "($XXM)" purse suffix
$ws.Range("B160").Value = "The Players Championship"
$ws.Range("C160").Value = "AT&T Pebble Beach Pro-Am"
$ws.Range("D160").Value = "Genesis Invitational"
$ws.Range("E160").Value = "Arnold Palmer Invitational"
$ws.Range("F160").Value = "RBC Heritage"
$ws.Range("G160").Value = "Memorial Tournament"
$ws.Range("H160").Value = "Travelers Championship"
$ws.Range("I160").Value = "FedEx St. Jude Championship"
$ws.Range("J160").Value = "BMW Championship"

# Replace the payout values in rows 161-170 (columns B-J) with a dash
# pending official confirmation
$ws.Range("B161:J170").Value = "-"
